$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.7
$ws.Range("H2").Value = -0.07563566285451317

$ws.Range("A3").Value = 1.3
$ws.Range("H3").Value = -0.04038692750915865
$ws.Range("I3").Value = 13.63636363636376

$ws.Range("A4").Value = 1.54
$ws.Range("H4").Value = -0.07830535807649208
$ws.Range("I4").Value = -10.90909090909079

$ws.Range("H5").Value = 0.07749399946415342
$ws.Range("I5").Value = -0.7407407407407398

$ws.Range("H6").Value = 0.02744298723289798
$ws.Range("I6").Value = 1.362637362637365

$ws.Range("H7").Value = 0.004971703331912686

$ws.Range("H8").Value = 0.003958081394145274

$ws.Range("A9").Value = 4.9
$ws.Range("H9").Value = -0.02732858725059259

$ws.Range("A10").Value = 5.34
$ws.Range("H10").Value = -0.009624417043498451
$ws.Range("I10").Value = 0.9756097560975602

$ws.Range("A11").Value = 5.5
$ws.Range("H11").Value = 0.009937416469076144
$ws.Range("I11").Value = 0.3041825095057037

$ws.Range("H12").Value = 0.0962391949236181
$ws.Range("I12").Value = -0.1204685710062451

$ws.Range("H13").Value = 0.1108251278555828
$ws.Range("I13").Value = 0.1674570243034973

$ws.Range("H14").Value = 0.1165009817314631
$ws.Range("I14").Value = 0.0795311845960653

$ws.Range("H15").Value = 0.1227010535482061
$ws.Range("I15").Value = -0.003870967741935516

$ws.Range("H16").Value = 0.128686327541306
$ws.Range("I16").Value = 0.1232153334637199

$ws.Range("H17").Value = 0.129310802540286
$ws.Range("I17").Value = -0.043124101581217

$ws.Range("H18").Value = 0.1314085491721864
$ws.Range("I18").Value = 0.07067137809187272

$ws.Range("H19").Value = 0.1234197296651822
$ws.Range("I19").Value = -0.3991945036721156

$ws.Range("H20").Value = 0.1251522166491153
$ws.Range("I20").Value = -0.0950413223140497

$ws.Range("H21").Value = 0.1248164347236817
